$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing "Code Review 1" scores in column B (was 25, now 16.67)
# and add matching scores in the new column C ("Code Review 1" second entry)
# for rows 4 through 9.
for ($row = 4; $row -le 9; $row++) {
    $ws.Cells.Item($row, 2).Value = 16.67
    $ws.Cells.Item($row, 3).Value = 16.67
    $ws.Cells.Item($row, 3).HorizontalAlignment = -4108
}

# Update the active selection to D9, matching the saved view state.
$ws.Range("D9").Select()
